$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: run_id=17, rss_url_id=1, date, response=200, item_count=9
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "2024-06-15 06:15:59"
$ws.Range("D18").Value = 200
$ws.Range("E18").Value = 9

# Row 19: run_id=18, rss_url_id=2, date, response=200, item_count=0
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "2024-06-15 06:15:59"
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = 0
